$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

$ws.Range("A6").Value = "Пивоварня"

$ws.Range("A7").Select()
